# 自动更新Excel文件 - 2026-01-14 23:12:40
# Daily refresh: decrement the "剩余" (remaining) day-count in column E for
# every shop row, except rows that are not in an active cycle (row 36, whose
# 剩余 already equals 总天 / has a malformed start date) and row 94, whose
# cycle expired (剩余 hit 1) and was renewed with a fresh start date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    if ($r -eq 36) {
        # Unchanged this cycle - still at full remaining days.
        continue
    }

    if ($r -eq 94) {
        # Cycle renewed: remaining days reset to 7, start date moved forward.
        $ws.Cells.Item($r, 5).Value2 = 7
        $ws.Cells.Item($r, 6).Value2 = 20260115
        continue
    }

    $remaining = $ws.Cells.Item($r, 5).Value2
    if ($remaining -ne $null) {
        $ws.Cells.Item($r, 5).Value2 = $remaining - 1
    }
}
